# Commit: "insight of analysis part" -
# "insight part changed from highest to lowest external debt."
#
# Slide 5 ("Insights/Deliverables of the analysis") contains two
# occurrences of the phrase "seems to have the highest external debt" -
# one under "Which country has the highest debt?" and one under
# "Which country has the least debt?". The second occurrence (the
# answer for the *least* debt question) should read "lowest" instead
# of "highest".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$shape = $slide.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

$searchText = " seems to have the highest external debt"
$replacementText = " seems to have the lowest external debt"

$occurrence = 0
$searchStart = 0
while ($true) {
    $found = $textRange.Find($searchText, $searchStart)
    if (-not $found) {
        break
    }
    $occurrence = $occurrence + 1
    if ($occurrence -eq 2) {
        $found.Text = $replacementText
        break
    }
    $searchStart = $found.Start + $found.Length
}
